$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 17, pushing existing rows 17:21 down to 18:22
$ws.Rows.Item(17).Insert()

# Populate the newly inserted row 17 with the new weekly data point
$ws.Cells.Item(17, 1).Value = 5
$ws.Cells.Item(17, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(17, 3).Value = "Maule"
$ws.Cells.Item(17, 4).Value = 44460
$ws.Cells.Item(17, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(17, 5).Value = 7
$ws.Cells.Item(17, 6).Value = 100112026
$ws.Cells.Item(17, 7).Value = "Haba"
$ws.Cells.Item(17, 8).Value = "Sin especificar"
$ws.Cells.Item(17, 9).Value = "Primera"
$ws.Cells.Item(17, 10).Value = 150
$ws.Cells.Item(17, 11).Value = 12000
$ws.Cells.Item(17, 12).Value = 12000
$ws.Cells.Item(17, 13).Value = 12000
$ws.Cells.Item(17, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(17, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(17, 16).Value = 480
$ws.Cells.Item(17, 17).Value = 25
$ws.Cells.Item(17, 18).Value = "Hortaliza"
